$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 46.727272
$ws.Range("I5").Value = 38.166668
$ws.Range("J5").Value = 57
$ws.Range("K5").Value = 38.166668
$ws.Range("L5").Value = 57
$ws.Range("M5").Value = 76.833332
$ws.Range("N5").Value = -287

$ws.Range("H40").Value = 1741.0869
$ws.Range("I40").Value = 1304.4
$ws.Range("J40").Value = 2077
$ws.Range("K40").Value = 1304.4
$ws.Range("L40").Value = 2077
$ws.Range("M40").Value = -1129.4
$ws.Range("N40").Value = -2427

$ws.Range("H88").Value = 11767889
$ws.Range("I88").Value = 4100.75
$ws.Range("J88").Value = 18490054
$ws.Range("K88").Value = 4100.75
$ws.Range("L88").Value = 18490054
$ws.Range("M88").Value = -3694.75
$ws.Range("N88").Value = -18490866

$ws.Range("H91").Value = 11767889
$ws.Range("I91").Value = 4100.75
$ws.Range("J91").Value = 18490054
$ws.Range("K91").Value = 4100.75
$ws.Range("L91").Value = 18490054
$ws.Range("M91").Value = -2696.75
$ws.Range("N91").Value = -18492862

$ws.Range("H111").Value = 1362.4166
$ws.Range("I111").Value = 1194.3334
$ws.Range("J111").Value = 1866.6666
$ws.Range("K111").Value = 3583.0002
$ws.Range("L111").Value = 5599.9998
$ws.Range("M111").Value = -516.0001999999999
$ws.Range("N111").Value = -11733.9998

$ws.Range("H131").Value = 1594.7222
$ws.Range("I131").Value = 541.6667
$ws.Range("J131").Value = 3700.8333
$ws.Range("K131").Value = 1625.0001
$ws.Range("L131").Value = 11102.4999
$ws.Range("M131").Value = 3414.9999
$ws.Range("N131").Value = -21182.4999

$ws.Range("H132").Value = 3646.5789
$ws.Range("I132").Value = 3775.7273
$ws.Range("K132").Value = 11327.1819
$ws.Range("M132").Value = -8797.1819

$ws.Range("H137").Value = 1772.3334
$ws.Range("I137").Value = 1754.0476
$ws.Range("J137").Value = 1900.3334
$ws.Range("K137").Value = 5262.142800000001
$ws.Range("L137").Value = 5701.0002
$ws.Range("M137").Value = -2712.142800000001
$ws.Range("N137").Value = -10801.0002

$ws.Range("H141").Value = 708.2222
$ws.Range("I141").Value = 600.6799999999999
$ws.Range("K141").Value = 1802.04
$ws.Range("M141").Value = 3377.96

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2683.2856
$ws.Range("I4").Value = 1463.8334
$ws.Range("K4").Value = 1463.8334
$ws.Range("M4").Value = -1347.8334

$ws.Range("H5").Value = 216.66667
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 275
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 275
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -499

$ws.Range("H26").Value = 2375
$ws.Range("I26").Value = 1750
$ws.Range("K26").Value = 1750
$ws.Range("M26").Value = -1420

$ws.Range("H41").Value = 1550
$ws.Range("I41").Value = 1550
$ws.Range("K41").Value = 1550
$ws.Range("M41").Value = -1136

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 216.66667
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 275
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 275
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -505

$ws.Range("H22").Value = 391
$ws.Range("I22").Value = 391
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 391
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -218
$ws.Range("N22").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()

$ws.Range("H80").Value = 287.33334
$ws.Range("I80").Value = 80
$ws.Range("K80").Value = 80
$ws.Range("M80").Value = 918

$ws.Range("H83").Value = 287.33334
$ws.Range("I83").Value = 80
$ws.Range("K83").Value = 400
$ws.Range("M83").Value = 4592

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6946650.5
$ws.Range("I31").Value = 10418525
$ws.Range("J31").Value = 2901.1667
$ws.Range("K31").Value = 10418525
$ws.Range("L31").Value = 2901.1667
$ws.Range("M31").Value = -10418230
$ws.Range("N31").Value = -3491.1667

$ws.Range("H34").Value = 6946650.5
$ws.Range("I34").Value = 10418525
$ws.Range("J34").Value = 2901.1667
$ws.Range("K34").Value = 10418525
$ws.Range("L34").Value = 2901.1667
$ws.Range("M34").Value = -10418323
$ws.Range("N34").Value = -3305.1667

$ws.Range("H39").Value = 3287.75
$ws.Range("I39").Value = 225.5
$ws.Range("J39").Value = 6350
$ws.Range("K39").Value = 225.5
$ws.Range("L39").Value = 6350
$ws.Range("M39").Value = 165.5
$ws.Range("N39").Value = -7132

$ws.Range("H49").Value = 3287.75
$ws.Range("I49").Value = 225.5
$ws.Range("J49").Value = 6350
$ws.Range("K49").Value = 225.5
$ws.Range("L49").Value = 6350
$ws.Range("M49").Value = -43.5
$ws.Range("N49").Value = -6714

$ws.Range("H132").Value = 33367.625
$ws.Range("I132").Value = 1548.6538
$ws.Range("K132").Value = 4645.9614
$ws.Range("M132").Value = -2115.9614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 236.81818
$ws.Range("I17").Value = 98.75
$ws.Range("J17").Value = 315.7143
$ws.Range("K17").Value = 296.25
$ws.Range("L17").Value = 947.1428999999999
$ws.Range("M17").Value = -127.25
$ws.Range("N17").Value = -1285.1429

$ws.Range("H80").Value = 1804.8334
$ws.Range("J80").Value = 1975.8
$ws.Range("L80").Value = 5927.4
$ws.Range("N80").Value = -7799.4

$ws.Range("H83").Value = 1804.8334
$ws.Range("J83").Value = 1975.8
$ws.Range("L83").Value = 17782.2
$ws.Range("N83").Value = -27142.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1536.2307
$ws.Range("I61").Value = 1243.5714
$ws.Range("J61").Value = 1877.6666
$ws.Range("K61").Value = 1243.5714
$ws.Range("L61").Value = 1877.6666
$ws.Range("M61").Value = -1041.5714
$ws.Range("N61").Value = -2281.6666

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H93").Value = 2382.5
$ws.Range("I93").Value = 2559.3333
$ws.Range("J93").Value = 1852
$ws.Range("K93").Value = 2559.3333
$ws.Range("L93").Value = 1852
$ws.Range("M93").Value = -1311.3333
$ws.Range("N93").Value = -4348

$ws.Range("H100").Value = 2338.8096
$ws.Range("I100").Value = 1901
$ws.Range("J100").Value = 3433.3333
$ws.Range("K100").Value = 1901
$ws.Range("L100").Value = 3433.3333
$ws.Range("M100").Value = -1360
$ws.Range("N100").Value = -4515.3333

$ws.Range("H113").Value = 1536.2307
$ws.Range("I113").Value = 1243.5714
$ws.Range("J113").Value = 1877.6666
$ws.Range("K113").Value = 1243.5714
$ws.Range("L113").Value = 1877.6666
$ws.Range("M113").Value = 926.4286
$ws.Range("N113").Value = -6217.6666

$ws.Range("H136").Value = 4725.0244
$ws.Range("I136").Value = 4388.6665
$ws.Range("J136").Value = 5642.364
$ws.Range("K136").Value = 13165.9995
$ws.Range("L136").Value = 16927.092
$ws.Range("M136").Value = -10615.9995
$ws.Range("N136").Value = -22027.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 201566.4
$ws.Range("I21").Value = 1407.5
$ws.Range("J21").Value = 335005.66
$ws.Range("K21").Value = 1407.5
$ws.Range("L21").Value = 335005.66
$ws.Range("M21").Value = -1172.5
$ws.Range("N21").Value = -335475.66

$ws.Range("H35").Value = 201566.4
$ws.Range("I35").Value = 1407.5
$ws.Range("J35").Value = 335005.66
$ws.Range("K35").Value = 1407.5
$ws.Range("L35").Value = 335005.66
$ws.Range("M35").Value = -1117.5
$ws.Range("N35").Value = -335585.66

$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002

$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 9000
$ws.Range("M83").Value = -4008

$ws.Range("H132").Value = 67424310
$ws.Range("J132").Value = 3238924
$ws.Range("L132").Value = 9716772
$ws.Range("N132").Value = -9721832
